$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-11 Thursday", "2024-01-12 Friday"),
    @("74×77=", "97×57="),
    @("22×25=", "27×27="),
    @("82×74=", "75×90="),
    @("44×64=", "48×55="),
    @("54×77=", "13×89="),
    @("26×54=", "95×62="),
    @("14×67=", "71×46="),
    @("50×49=", "68×31="),
    @("30×17=", "94×99="),
    @("82×44=", "15×87="),
    @("66×46=", "96×66="),
    @("26×43=", "27×35="),
    @("21×70=", "11×86="),
    @("84×42=", "48×61="),
    @("17×56=", "69×93="),
    @("93×90=", "42×64="),
    @("39×64=", "81×82="),
    @("89×85=", "39×44="),
    @("64×77=", "98×63="),
    @("53×98=", "19×17="),
    @("54×98=", "75×17="),
    @("90×27=", "88×22="),
    @("12×80=", "23×19="),
    @("89×97=", "33×73="),
    @("16×90=", "81×65=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
